$wb = $excel.ActiveWorkbook

# Update timestamps on "data" sheet
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Range("F2").Value = "2021-10-05 14:20:09.356753"
$dataSheet.Range("F3").Value = "2021-10-05 14:20:09.356760"

# Add new "metadata" sheet after "data"
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Match the page margins used on the "data" sheet
$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

# Header row values
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row values
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Erythropoietic protoporphyria, mild variant"
$metaSheet.Range("C2").Value = 91
# "1.2" must stay textual (not be coerced into the number 1.2)
$metaSheet.Range("D2").Value = "'1.2"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2017-11-05T02:37:20.011517Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:09.353556"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/91/?format=json"

# Reuse the same "header" cell style used on the "data" sheet (bold, centered, bordered)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
